$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 2253.6191
$ws.Range("J106").Value = 3192.2
$ws.Range("L106").Value = 3192.2
$ws.Range("N106").Value = -4454.2
$ws.Range("H107").Value = 842.1667
$ws.Range("I107").Value = 894.3
$ws.Range("J107").Value = 581.5
$ws.Range("K107").Value = 894.3
$ws.Range("L107").Value = 581.5
$ws.Range("M107").Value = 1025.7
$ws.Range("N107").Value = -4421.5
$ws.Range("H125").Value = 1247.0769
$ws.Range("I125").Value = 711
$ws.Range("J125").Value = 1344.5454
$ws.Range("K125").Value = 6399
$ws.Range("L125").Value = 12100.9086
$ws.Range("M125").Value = -3939
$ws.Range("N125").Value = -17020.9086
$ws.Range("H137").Value = 17428.016
$ws.Range("I137").Value = 1408.641
$ws.Range("J137").Value = 43459.5
$ws.Range("K137").Value = 4225.923000000001
$ws.Range("L137").Value = 130378.5
$ws.Range("M137").Value = -1675.923000000001
$ws.Range("N137").Value = -135478.5

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15450.467
$ws.Range("I32").Value = 17060.37
$ws.Range("J32").Value = 4986.1
$ws.Range("K32").Value = 17060.37
$ws.Range("L32").Value = 4986.1
$ws.Range("M32").Value = -16773.37
$ws.Range("N32").Value = -5560.1
$ws.Range("H37").Value = 29990
$ws.Range("J37").Value = 29990
$ws.Range("L37").Value = 29990
$ws.Range("N37").Value = -30536
$ws.Range("H41").Value = 3503
$ws.Range("I41").Value = 3503
$ws.Range("K41").Value = 3503
$ws.Range("M41").Value = -3089
$ws.Range("H44").Value = 29237.5
$ws.Range("J44").Value = 29237.5
$ws.Range("L44").Value = 29237.5
$ws.Range("N44").Value = -30213.5
$ws.Range("H61").Value = 2819.2856
$ws.Range("I61").Value = 2099.8965
$ws.Range("K61").Value = 2099.8965
$ws.Range("M61").Value = -1887.8965
$ws.Range("H102").Value = 5019
$ws.Range("I102").Value = 3730.4
$ws.Range("K102").Value = 3730.4
$ws.Range("M102").Value = -2108.4
$ws.Range("H132").Value = 19876.75
$ws.Range("I132").Value = 1465.4117
$ws.Range("J132").Value = 48330.637
$ws.Range("K132").Value = 4396.2351
$ws.Range("L132").Value = 144991.911
$ws.Range("M132").Value = -1866.2351
$ws.Range("N132").Value = -150051.911
$ws.Range("H136").Value = 2819.2856
$ws.Range("I136").Value = 2099.8965
$ws.Range("K136").Value = 6299.689499999999
$ws.Range("M136").Value = -3749.689499999999

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1257.1875
$ws.Range("I20").Value = 862.6923
$ws.Range("K20").Value = 862.6923
$ws.Range("M20").Value = -615.6923

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 38.285713
$ws.Range("I7").Value = 36.333332
$ws.Range("K7").Value = 36.333332
$ws.Range("M7").Value = 76.666668
$ws.Range("H31").Value = 9536.34
$ws.Range("I31").Value = 17550.291
$ws.Range("J31").Value = 2904.1035
$ws.Range("K31").Value = 17550.291
$ws.Range("L31").Value = 2904.1035
$ws.Range("M31").Value = -17255.291
$ws.Range("N31").Value = -3494.1035
$ws.Range("H34").Value = 9536.34
$ws.Range("I34").Value = 17550.291
$ws.Range("J34").Value = 2904.1035
$ws.Range("K34").Value = 17550.291
$ws.Range("L34").Value = 2904.1035
$ws.Range("M34").Value = -17348.291
$ws.Range("N34").Value = -3308.1035
$ws.Range("H58").Value = 18429.793
$ws.Range("I58").Value = 1224.5927
$ws.Range("J58").Value = 250700
$ws.Range("K58").Value = 1224.5927
$ws.Range("L58").Value = 250700
$ws.Range("M58").Value = -1021.5927
$ws.Range("N58").Value = -251106
$ws.Range("H132").Value = 16251.081
$ws.Range("I132").Value = 17734.355
$ws.Range("K132").Value = 53203.065
$ws.Range("M132").Value = -50673.065
$ws.Range("H136").Value = 18429.793
$ws.Range("I136").Value = 1224.5927
$ws.Range("J136").Value = 250700
$ws.Range("K136").Value = 3673.7781
$ws.Range("L136").Value = 752100
$ws.Range("M136").Value = -1123.7781
$ws.Range("N136").Value = -757200

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 523.0909
$ws.Range("I17").Value = 200
$ws.Range("J17").Value = 644.25
$ws.Range("K17").Value = 600
$ws.Range("L17").Value = 1932.75
$ws.Range("M17").Value = -431
$ws.Range("N17").Value = -2270.75
$ws.Range("H34").Value = 613.3333
$ws.Range("J34").Value = 850
$ws.Range("L34").Value = 2550
$ws.Range("N34").Value = -2718
$ws.Range("H93").Value = 3000
$ws.Range("J93").Value = 3500
$ws.Range("L93").Value = 10500
$ws.Range("N93").Value = -14244
$ws.Range("H103").Value = 988.25
$ws.Range("I103").Value = 575
$ws.Range("J103").Value = 1401.5
$ws.Range("K103").Value = 1725
$ws.Range("L103").Value = 4204.5
$ws.Range("M103").Value = -846
$ws.Range("N103").Value = -5962.5
$ws.Range("H122").Value = 1091
$ws.Range("I122").Value = 361.0909
$ws.Range("J122").Value = 1664.5
$ws.Range("K122").Value = 3249.8181
$ws.Range("L122").Value = 14980.5
$ws.Range("M122").Value = -799.8181
$ws.Range("N122").Value = -19880.5
$ws.Range("H131").Value = 814.73
$ws.Range("I131").Value = 603.75
$ws.Range("J131").Value = 833.0761
$ws.Range("K131").Value = 1811.25
$ws.Range("L131").Value = 2499.2283
$ws.Range("M131").Value = 3228.75
$ws.Range("N131").Value = -12579.2283

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 3000000
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("H8").Value = 3000000
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1875.5714
$ws.Range("J16").Value = 1863.6
$ws.Range("L16").Value = 1863.6
$ws.Range("N16").Value = -2203.6
$ws.Range("H40").Value = 64230.168
$ws.Range("I40").Value = 94271.914
$ws.Range("J40").Value = 4146.6665
$ws.Range("K40").Value = 94271.914
$ws.Range("L40").Value = 4146.6665
$ws.Range("M40").Value = -94135.914
$ws.Range("N40").Value = -4418.6665
$ws.Range("H55").Value = 784.55554
$ws.Range("I55").Value = 974.5454999999999
$ws.Range("J55").Value = 486
$ws.Range("K55").Value = 974.5454999999999
$ws.Range("L55").Value = 486
$ws.Range("M55").Value = -801.5454999999999
$ws.Range("N55").Value = -832
$ws.Range("H82").Value = 2429.4
$ws.Range("J82").Value = 1276.5
$ws.Range("L82").Value = 1276.5
$ws.Range("N82").Value = -1998.5
$ws.Range("H85").Value = 2429.4
$ws.Range("J85").Value = 1276.5
$ws.Range("L85").Value = 1276.5
$ws.Range("N85").Value = -3772.5
$ws.Range("H93").Value = 2284.3635
$ws.Range("I93").Value = 2452.4
$ws.Range("K93").Value = 2452.4
$ws.Range("M93").Value = -1204.4
$ws.Range("H122").Value = 3786.2856
$ws.Range("I122").Value = 3300.4285
$ws.Range("J122").Value = 4272.143
$ws.Range("K122").Value = 9901.2855
$ws.Range("L122").Value = 12816.429
$ws.Range("M122").Value = -7451.2855
$ws.Range("N122").Value = -17716.429
$ws.Range("H132").Value = 1867.9286
$ws.Range("I132").Value = 1100.2858
$ws.Range("K132").Value = 3300.8574
$ws.Range("M132").Value = -770.8574000000003
$ws.Range("H136").Value = 16959.303
$ws.Range("I136").Value = 25048.762
$ws.Range("J136").Value = 2802.75
$ws.Range("K136").Value = 75146.28599999999
$ws.Range("L136").Value = 8408.25
$ws.Range("M136").Value = -72596.28599999999
$ws.Range("N136").Value = -13508.25

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 25000
$ws.Range("I70").Value = 25000
$ws.Range("J70").Value = 25000
$ws.Range("K70").Value = 25000
$ws.Range("L70").Value = 25000
$ws.Range("M70").Value = -24685
$ws.Range("N70").Value = -25630
$ws.Range("H73").Value = 25000
$ws.Range("I73").Value = 25000
$ws.Range("J73").Value = 25000
$ws.Range("K73").Value = 25000
$ws.Range("L73").Value = 25000
$ws.Range("M73").Value = -23908
$ws.Range("N73").Value = -27184
$ws.Range("H113").Value = 1352082.6
$ws.Range("I113").Value = 818.94116
$ws.Range("J113").Value = 9009244
$ws.Range("K113").Value = 2456.82348
$ws.Range("L113").Value = 27027732
$ws.Range("M113").Value = -286.82348
$ws.Range("N113").Value = -27032072
$ws.Range("H122").Value = 1604.5416
$ws.Range("I122").Value = 1457.5
$ws.Range("J122").Value = 3222
$ws.Range("K122").Value = 4372.5
$ws.Range("L122").Value = 9666
$ws.Range("M122").Value = -1922.5
$ws.Range("N122").Value = -14566
$ws.Range("H126").Value = 1653.6666
$ws.Range("I126").Value = 1369
$ws.Range("J126").Value = 2650
$ws.Range("K126").Value = 4107
$ws.Range("L126").Value = 7950
$ws.Range("M126").Value = -1637
$ws.Range("N126").Value = -12890
$ws.Range("H132").Value = 1708.2683
$ws.Range("I132").Value = 1449.0322
$ws.Range("J132").Value = 2511.9
$ws.Range("K132").Value = 4347.096600000001
$ws.Range("L132").Value = 7535.700000000001
$ws.Range("M132").Value = -1817.096600000001
$ws.Range("N132").Value = -12595.7
$ws.Range("H136").Value = 1271.138
$ws.Range("I136").Value = 798.2857
$ws.Range("J136").Value = 2512.375
$ws.Range("K136").Value = 2394.8571
$ws.Range("L136").Value = 7537.125
$ws.Range("M136").Value = 155.1428999999998
$ws.Range("N136").Value = -12637.125
